$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use Text number format while writing, to keep numeric-looking values (e.g. "336.16")
# stored as text instead of being parsed into floating point numbers; then restore
# the default "Normal" style so no stray style attribute is left on the cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "29.981.74"
$ws.Range("E2").Value = "  +1.71%  "
$ws.Range("D3").Value = "1.941.86"
$ws.Range("E3").Value = "  +1.15%  "
$ws.Range("D4").Value = "1.008"
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "336.16"
$ws.Range("E5").Value = "  +3.08%  "
$ws.Range("E7").Value = "  +0.30%  "
$ws.Range("D8").Value = "0.4143"
$ws.Range("E8").Value = "  +1.37%  "
$ws.Range("D9").Value = "0.08215"
$ws.Range("E9").Value = "  -0.37%  "
$ws.Range("D10").Value = "1.018"
$ws.Range("E10").Value = "  -0.59%  "
$ws.Range("D11").Value = "23.94"
$ws.Range("E11").Value = "  +1.83%  "
$ws.Range("D12").Value = "1.956.47"
$ws.Range("E12").Value = "  +1.51%  "
$ws.Range("D13").Value = "6.110"
$ws.Range("E13").Value = "  +1.01%  "
$ws.Range("D14").Value = "7.322"
$ws.Range("E14").Value = "  +1.09%  "
$ws.Range("D15").Value = "91.45"
$ws.Range("E15").Value = "  +0.16%  "
$ws.Range("D16").Value = "0.06871"
$ws.Range("E16").Value = "  +1.06%  "
$ws.Range("D17").Value = "1.010"
$ws.Range("E17").Value = "  +0.12%  "
$ws.Range("D18").Value = "0.00001041"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("E19").Value = "  +0.67%  "
$ws.Range("D21").Value = "29.979.16"
$ws.Range("E21").Value = "  +1.56%  "
$ws.Range("D22").Value = "5.653"
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("E23").Value = "  +1.20%  "
$ws.Range("D24").Value = "2.198"
$ws.Range("E24").Value = "  +0.19%  "
$ws.Range("D25").Value = "2.199.57"
$ws.Range("E25").Value = "  +1.44%  "
$ws.Range("D26").Value = "6.693"
$ws.Range("E26").Value = "  +1.01%  "
$ws.Range("D27").Value = "157.05"
$ws.Range("E27").Value = "  +0.19%  "
$ws.Range("E28").Value = "  +0.05%  "
$ws.Range("D29").Value = "2.109"
$ws.Range("E29").Value = "  -0.48%  "
$ws.Range("D30").Value = "121.47"
$ws.Range("E30").Value = "  +0.81%  "
$ws.Range("D31").Value = "1.018"
$ws.Range("E31").Value = "  -0.66%  "
$ws.Range("D32").Value = "0.09643"
$ws.Range("E32").Value = "  +0.68%  "
$ws.Range("D33").Value = "5.626"
$ws.Range("E33").Value = "  +1.75%  "
$ws.Range("D34").Value = "1.424"
$ws.Range("E34").Value = "  +2.91%  "
$ws.Range("D35").Value = "3.556"
$ws.Range("E35").Value = "  -0.18%  "
$ws.Range("D36").Value = "0.06565"
$ws.Range("E36").Value = "  +7.01%  "
$ws.Range("D37").Value = "0.02294"
$ws.Range("E37").Value = "  +0.44%  "
$ws.Range("D38").Value = "1.220"
$ws.Range("E38").Value = "  +3.50%  "
$ws.Range("D39").Value = "0.5986"
$ws.Range("E39").Value = "  -0.06%  "
$ws.Range("D40").Value = "8.019"
$ws.Range("E40").Value = "  -0.37%  "
$ws.Range("E41").Value = "  -0.53%  "
$ws.Range("D42").Value = "2.535"
$ws.Range("E42").Value = "  +5.54%  "
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("E44").Value = "  -0.12%  "
$ws.Range("D45").Value = "12.45"
$ws.Range("E45").Value = "  +0.14%  "
$ws.Range("D46").Value = "0.07524"
$ws.Range("E46").Value = "  -1.09%  "
$ws.Range("D47").Value = "0.5577"
$ws.Range("E47").Value = "  -0.09%  "
$ws.Range("D48").Value = "1.990"
$ws.Range("E48").Value = "  +1.52%  "
$ws.Range("D49").Value = "117.75"
$ws.Range("E49").Value = "  +0.01%  "
$ws.Range("D50").Value = "73.01"
$ws.Range("E50").Value = "  +0.34%  "
$ws.Range("D51").Value = "2.425"
$ws.Range("E51").Value = "  -0.23%  "

$ws.Range("D2:D51").Style = "Normal"
